$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.751.84"
$ws.Range("E2").Value = "  -0.11%  "
$ws.Range("D3").Value = "2.293.76"
$ws.Range("E3").Value = "  +3.58%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "269.89"
$ws.Range("E5").Value = "  +2.57%  "
$ws.Range("D6").Value = "93.71"
$ws.Range("E6").Value = "  +8.43%  "
$ws.Range("D7").Value = "0.633"
$ws.Range("E7").Value = "  +1.85%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.619"
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").Value = "45.37"
$ws.Range("E10").Value = "  -2.16%  "
$ws.Range("D11").Value = "0.0936"
$ws.Range("E11").Value = "  +1.69%  "
$ws.Range("D12").Value = "8.06"
$ws.Range("E12").Value = "  +5.95%  "
$ws.Range("E13").Value = "  +0.38%  "
$ws.Range("D14").Value = "2.636.89"
$ws.Range("E14").Value = "  +3.42%  "
$ws.Range("D15").Value = "15.24"
$ws.Range("E15").Value = "  +4.27%  "
$ws.Range("D16").Value = "0.852"
$ws.Range("E16").Value = "  +8.26%  "
$ws.Range("D17").Value = "2.290.08"
$ws.Range("E17").Value = "  +1.25%  "
$ws.Range("D18").Value = "43.702.45"
$ws.Range("E18").Value = "  -0.29%  "
$ws.Range("E19").Value = "  +1.31%  "
$ws.Range("D20").Value = "6.26"
$ws.Range("E20").Value = "  +4.48%  "
$ws.Range("D21").Value = "71.12"
$ws.Range("E21").Value = "  +1.45%  "
$ws.Range("D22").Value = "2.29"
$ws.Range("E22").Value = "  -4.00%  "
$ws.Range("D23").Value = "236.36"
$ws.Range("E23").Value = "  +1.82%  "
$ws.Range("D24").Value = "9.68"
$ws.Range("E24").Value = "  +7.05%  "
$ws.Range("D25").Value = "1.00"
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "2.52"
$ws.Range("E26").Value = "  +10.74%  "
$ws.Range("D27").Value = "11.23"
$ws.Range("E27").Value = "  +3.88%  "
$ws.Range("D28").Value = "3.40"
$ws.Range("E28").Value = "  -3.34%  "
$ws.Range("D29").Value = "39.26"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").Value = "2.26"
$ws.Range("E30").Value = "  +0.13%  "
$ws.Range("D31").Value = "22.24"
$ws.Range("E31").Value = "  +8.37%  "
$ws.Range("D32").Value = "173.53"
$ws.Range("E32").Value = "  -0.84%  "
$ws.Range("D33").Value = "0.0886"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "5.54"
$ws.Range("E34").Value = "  +1.96%  "
$ws.Range("E35").Value = "  +2.01%  "
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").Value = "4.53"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").Value = "0.0349"
$ws.Range("E38").Value = "  -3.21%  "
$ws.Range("D39").Value = "3.39"
$ws.Range("E39").Value = "  +3.63%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").Value = "2.32"
$ws.Range("E40").Value = "  +10.41%  "
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "0.233"
$ws.Range("E41").Value = "  +13.52%  "
$ws.Range("D42").Value = "12.27"
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").Value = "1.32"
$ws.Range("E43").Value = "  +16.71%  "
$ws.Range("D44").Value = "5.45"
$ws.Range("E44").Value = "  -2.00%  "
$ws.Range("D45").Value = "61.24"
$ws.Range("E45").Value = "  -5.16%  "
$ws.Range("E46").Value = "  +5.90%  "
$ws.Range("D47").Value = "0.102"
$ws.Range("E47").Value = "  +3.33%  "
$ws.Range("D48").Value = "100.03"
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("D49").Value = "1.19"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("B50").Value = "WOONetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D50").Value = "0.432"
$ws.Range("E50").Value = "  -3.48%  "
$ws.Range("B51").Value = "RocketPoolETH"
$ws.Range("C51").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D51").Value = "2.516.83"
$ws.Range("E51").Value = "  +3.39%  "
